$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J1").Value = "State"
$ws.Range("D4").Value = "Lines"

# Copy header-row formatting (bold + medium box border) onto D8
$ws.Range("A1").Copy()
$ws.Range("D8").PasteSpecial(-4122)  # xlPasteFormats

# Copy row-2 formatting (thin border, no top) onto D9
$ws.Range("A2").Copy()
$ws.Range("D9").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("D8").Value = "Line"
$ws.Range("D9").Value = "cells"

$ws.Rows.Item(7).RowHeight = 15.75
$ws.Rows.Item(8).RowHeight = 15.75

$ws.Range("C3").Select()
